$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "56.248.71"
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("D3").Value = "2.974.11"
$ws.Range("E3").Value = "  -5.00%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'495.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.91%  "
$ws.Range("D6").Value = "'135.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "2.972.46"
$ws.Range("E8").Value = "  -5.15%  "
$ws.Range("D9").Value = "'0.426"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "'7.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.03%  "
$ws.Range("D12").Value = "'0.351"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.85%  "
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "3.488.61"
$ws.Range("E14").Value = "  -4.69%  "
$ws.Range("D15").Value = "'24.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "56.308.46"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "2.982.92"
$ws.Range("D18").Value = "'0.0000145"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.37%  "
$ws.Range("D19").Value = "'5.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'12.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.42%  "
$ws.Range("D21").Value = "'7.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "'323.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.94%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'0.462"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.85%  "
$ws.Range("D25").Value = "'61.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.03%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "'0.163"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "0.0₃0890"
$ws.Range("E28").Value = "  -6.39%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'6.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").Value = "'6.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").Value = "'1.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").Value = "'1.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.81%  "
$ws.Range("D34").Value = "'19.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.08%  "
$ws.Range("D35").Value = "'154.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("D36").Value = "'4.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.12%  "
$ws.Range("D37").Value = "'1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("D38").Value = "'5.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.96%  "
$ws.Range("D39").Value = "'0.0666"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("D40").Value = "'23.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").Value = "3.007.66"
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("D42").Value = "'37.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.12%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.67%  "
$ws.Range("D45").Value = "'1.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "'0.633"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.96%  "
$ws.Range("D47").Value = "2.207.54"
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("E48").Value = "  -8.78%  "
$ws.Range("D49").Value = "'1.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.15%  "
$ws.Range("D50").Value = "'0.0236"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'19.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.77%  "
